$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new bad-GC-run entry (duplicate of the 20OCT22FID_TC_09 run) as row 20.
$ws.Range("A20").Value = "data/GC_Data/all_data/20OCT22FID_TC_09.RES"
$ws.Range("A20").Font.Color = 0

# Update the current selection to D20 (matches the author's last on-screen selection).
[void]$ws.Range("D20").Select()
